# Add a "Skills" worksheet to the Profile/Sign Up workbook (between "Profile"
# and "Sign Up"), containing a small SkillName/C# table, and make
# "Credentials" the active tab again.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after "Profile" so the tab order becomes
# Credentials, Profile, Skills, Sign Up.
$profileSheet = $wb.Worksheets.Item("Profile")
$skills = $wb.Worksheets.Add($null, $profileSheet)
$skills.Name = "Skills"

# Content: a one-column table with a header and a single skill.
$skills.Range("A1").Value = "SkillName"
$skills.Range("A2").Value = "C#"

# Match the authored column width for column A.
$skills.Columns.Item(1).ColumnWidth = 18.14

# Restore "Credentials" as the selected/active tab (it was the active tab
# before this edit; "Sign Up" had been the active one and loses that flag).
$credentials = $wb.Worksheets.Item("Credentials")
$credentials.Activate()
$credentials.Range("B3").Select()
